$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that was left after "Fall 2021".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Fill in the "Department Chair" content control (tag "DC") with the
#    chair's name, replacing the "Click here to enter text." placeholder.
for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    if ($cc.Tag -eq "DC" -and $cc.ID -eq "-1339623062") {
        $cc.Range.Text = "Brady Nielsen"
        break
    }
}

# 3. Fill in the "Dean" content control (tag "Dean", the one belonging to
#    the chair/dean signature block) with the dean's name, replacing the
#    "Click here to enter text." placeholder. (There is another unrelated
#    content control elsewhere in the document that is also tagged "Dean"
#    and must be left untouched.)
for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    if ($cc.Tag -eq "Dean" -and $cc.ID -eq "-257067089") {
        $cc.Range.Text = "Sarah Martin"
        break
    }
}
